$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.825549
$ws.Range("H2").Value = 5.476647
$ws.Range("I2").Value = 0.04696949406168958
$ws.Range("J2").Value = 0.04696949406168958
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.42119933333333
$ws.Range("N2").Value = 73.263598
$ws.Range("O2").Value = 0.4086816635579248
$ws.Range("P2").Value = 0.4086816635579248
$ws.Range("Q2").Value = 44.58209602176733
$ws.Range("R2").Value = 401.238864195906
$ws.Range("S2").Value = 0.01919557096960537
$ws.Range("T2").Value = 0.01919557096960537
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.825549
$ws.Range("H3").Value = 5.476647
$ws.Range("I3").Value = 0.04696949406168958
$ws.Range("J3").Value = 0.04696949406168958
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 33.48129
$ws.Range("N3").Value = 100.44387
$ws.Range("O3").Value = 0.5602996441124273
$ws.Range("P3").Value = 0.5602996441124273
$ws.Range("Q3").Value = 61.12173547821
$ws.Range("R3").Value = 550.09561930389
$ws.Range("S3").Value = 0.02631699080690544
$ws.Range("T3").Value = 0.02631699080690544
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.825549
$ws.Range("H4").Value = 5.476647
$ws.Range("I4").Value = 0.04696949406168958
$ws.Range("J4").Value = 0.04696949406168958
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.853554333333333
$ws.Range("N4").Value = 5.560663
$ws.Range("O4").Value = 0.03101869232964781
$ws.Range("P4").Value = 0.03101869232964781
$ws.Range("Q4").Value = 3.383754259662333
$ws.Range("R4").Value = 30.453788336961
$ws.Range("S4").Value = 0.001456932285178769
$ws.Range("T4").Value = 0.001456932285178769
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.57737633333333
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.6066215903701957
$ws.Range("J5").Value = 0.6066215903701957
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 24.42119933333333
$ws.Range("N5").Value = 73.263598
$ws.Range("O5").Value = 0.4086816635579248
$ws.Range("P5").Value = 0.4086816635579248
$ws.Range("Q5").Value = 575.7878071933491
$ws.Range("R5").Value = 5182.090264740141
$ws.Range("S5").Value = 0.2479151207026456
$ws.Range("T5").Value = 0.2479151207026456
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.57737633333333
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.6066215903701957
$ws.Range("J6").Value = 0.6066215903701957
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 33.48129
$ws.Range("N6").Value = 100.44387
$ws.Range("O6").Value = 0.5602996441124273
$ws.Range("P6").Value = 0.5602996441124273
$ws.Range("Q6").Value = 789.40097445547
$ws.Range("R6").Value = 7104.608770099229
$ws.Range("S6").Value = 0.3398898611953353
$ws.Range("T6").Value = 0.3398898611953353
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.57737633333333
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.6066215903701957
$ws.Range("J7").Value = 0.6066215903701957
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.853554333333333
$ws.Range("N7").Value = 5.560663
$ws.Range("O7").Value = 0.03101869232964781
$ws.Range("P7").Value = 0.03101869232964781
$ws.Range("Q7").Value = 43.70194807128077
$ws.Range("R7").Value = 393.3175326415269
$ws.Range("S7").Value = 0.01881660847221475
$ws.Range("T7").Value = 0.01881660847221475
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 13.46376966666667
$ws.Range("H8").Value = 40.391309
$ws.Range("I8").Value = 0.3464089155681148
$ws.Range("J8").Value = 0.3464089155681148
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 24.42119933333333
$ws.Range("N8").Value = 73.263598
$ws.Range("O8").Value = 0.4086816635579248
$ws.Range("P8").Value = 0.4086816635579248
$ws.Range("Q8").Value = 328.8014028077536
$ws.Range("R8").Value = 2959.212625269782
$ws.Range("S8").Value = 0.1415709718856739
$ws.Range("T8").Value = 0.1415709718856739
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 13.46376966666667
$ws.Range("H9").Value = 40.391309
$ws.Range("I9").Value = 0.3464089155681148
$ws.Range("J9").Value = 0.3464089155681148
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 33.48129
$ws.Range("N9").Value = 100.44387
$ws.Range("O9").Value = 0.5602996441124273
$ws.Range("P9").Value = 0.5602996441124273
$ws.Range("Q9").Value = 450.78437670287
$ws.Range("R9").Value = 4057.05939032583
$ws.Range("S9").Value = 0.1940927921101866
$ws.Range("T9").Value = 0.1940927921101866
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 13.46376966666667
$ws.Range("H10").Value = 40.391309
$ws.Range("I10").Value = 0.3464089155681148
$ws.Range("J10").Value = 0.3464089155681148
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.853554333333333
$ws.Range("N10").Value = 5.560663
$ws.Range("O10").Value = 0.03101869232964781
$ws.Range("P10").Value = 0.03101869232964781
$ws.Range("Q10").Value = 24.95582860865189
$ws.Range("R10").Value = 224.602457477867
$ws.Range("S10").Value = 0.0107451515722543
$ws.Range("T10").Value = 0.0107451515722543
